$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-26 Monday", "2026-01-27 Tuesday"),
    @("343÷7=", "366÷9="),
    @("685÷6=", "421÷3="),
    @("985÷5=", "222÷3="),
    @("192÷2=", "799÷4="),
    @("995÷2=", "570÷5="),
    @("575÷2=", "182÷9="),
    @("404÷8=", "567÷4="),
    @("462÷8=", "665÷6="),
    @("997÷8=", "292÷9="),
    @("459÷5=", "943÷4="),
    @("658÷7=", "841÷3="),
    @("644÷8=", "757÷8="),
    @("545÷6=", "307÷8="),
    @("948÷2=", "556÷3="),
    @("172÷5=", "280÷2="),
    @("920÷7=", "236÷2="),
    @("391÷6=", "561÷6="),
    @("315÷4=", "786÷6="),
    @("533÷6=", "203÷5="),
    @("524÷2=", "595÷3="),
    @("925÷9=", "387÷6="),
    @("747÷2=", "575÷4="),
    @("790÷5=", "422÷5="),
    @("292÷2=", "613÷3="),
    @("101÷4=", "268÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
